# Add a new "Technology" question row to the "questions" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("questions")

$ws.Range("A6").Value = "Technology"
$ws.Range("B6").Value = "What are all the technologies that are used currently in Bain?"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "technology list"

# Move the selection down to the next empty row, matching where the
# cursor ends up after the new row is entered.
$ws.Range("A7").Select()
